# Daily auto-push update: insert one new timestamped ranking sample.
#
# The new reading (2026/02/28, 土, hour 13, ranking 201) belongs right
# after the existing 2026/02/28 07:00 row (row 901) and before the
# 2026/12/29 row (old row 902). Inserting a whole row there shifts every
# subsequent row down by one (old row 902 -> new row 903, ..., old row
# 943 -> new row 944), which matches the diff exactly and grows the used
# range from A1:D943 to A1:D944.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 902..943 down to 903..944 and free up row 902 for the new entry.
$ws.Rows.Item(902).Insert()

# Column A holds dates as literal text (e.g. "2026/12/29"), not real Excel
# date serials. Setting NumberFormat to Text ("@") before assigning the
# string stops Excel's automatic "looks like a date" conversion, so the
# cell keeps the literal string "2026/02/28" exactly like its neighbours.
$ws.Cells.Item(902, 1).NumberFormat = "@"
$ws.Cells.Item(902, 1).Value = "2026/02/28"
$ws.Cells.Item(902, 2).Value = "土"
$ws.Cells.Item(902, 3).Value = 13
$ws.Cells.Item(902, 4).Value = 201

# Re-pull the cell's number format from its neighbour (row 901, which is
# still plain/General like every other data row) so A902 ends up styled
# identically to the rest of the column instead of staying tagged as a
# distinct "Text" format cell.
$ws.Cells.Item(901, 1).Copy()
$ws.Cells.Item(902, 1).PasteSpecial(-4122)  # xlPasteFormats

Write-Output ("Inserted row 902: " + $ws.Cells.Item(902,1).Value2 + ", " + `
    $ws.Cells.Item(902,2).Value2 + ", " + $ws.Cells.Item(902,3).Value2 + `
    ", " + $ws.Cells.Item(902,4).Value2)
Write-Output ("New used range: " + $ws.UsedRange.Address())
